# Applies the workbook edit described in the commit:
# "Added asymmetry index plot to polar_loc_speed_motile.m"
# New experimental runs (2022-08-12 "bad?" run and 2022-08-25 run) of strain
# 1633 mNG_FimW FimX_mScI cpdA- pch- are recorded on the "Input" sheet and on
# the "FimWX Jenal double" data sheet, and the 2022-07-27 run on the data
# sheet is flagged "bad?" in column D.

$wb = $excel.ActiveWorkbook

$strain = "1633 mNG_FimW FimX_mScI cpdA- pch-"
$interval2h = "5s interval-2h37"
$interval1p5h = "5s interval-1p5h37"

# ---------------------------------------------------------------------------
# "FimWX Jenal double" sheet: flag the 2022-07-27 run as "bad?" first, so the
# new shared strings are introduced in the same order as in the saved file
# ("bad?" before "5s interval-1p5h37").
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("FimWX Jenal double")
$wsData.Range("D10").Value = "bad?"

# ---------------------------------------------------------------------------
# "Input" sheet: replace the 3 existing rows with the 4 rows that now drive
# the plotting script, taken from the "FimWX Jenal double" data (rows 11-14).
# ---------------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")

$wsInput.Range("A1").Style = "Normal"
$wsInput.Range("A1").Value = $strain
$wsInput.Range("B1").Value = 20220728
$wsInput.Range("C1").Style = "Normal"
$wsInput.Range("C1").Value = $interval2h

$wsInput.Range("A2").Style = "Normal"
$wsInput.Range("A2").Value = $strain
$wsInput.Range("B2").Value = 20220729
$wsInput.Range("C2").Style = "Normal"
$wsInput.Range("C2").Value = $interval2h

$wsInput.Range("A3").Value = $strain
$wsInput.Range("B3").Style = "Normal"
$wsInput.Range("B3").Value = 20220812
$wsInput.Range("C3").Value = $interval1p5h

$wsInput.Range("A4").Value = $strain
$wsInput.Range("B4").Value = 20220825
$wsInput.Range("C4").Value = $interval2h
$wsInput.Range("A4").Style = "Normal"
$wsInput.Range("C4").Style = "Normal"

$wsInput.Range("C3").Select() | Out-Null

# ---------------------------------------------------------------------------
# "FimWX Jenal double" sheet: append the two new runs (2022-08-12, 2022-08-25).
# ---------------------------------------------------------------------------
$wsData.Range("A13").Value = $strain
$wsData.Range("B13").Value = 20220812
$wsData.Range("C13").Value = $interval1p5h
$wsData.Range("A13").Style = "Normal"
$wsData.Range("C13").Style = "Normal"

$wsData.Range("A14").Value = $strain
$wsData.Range("B14").Value = 20220825
$wsData.Range("C14").Value = $interval2h
$wsData.Range("A14").Style = "Normal"
$wsData.Range("C14").Style = "Normal"

$wsData.Range("C13").Select() | Out-Null
$wsData.Activate() | Out-Null
